$d = $word.ActiveDocument

# Locate the target paragraph containing the original sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd("`r")
    if ($txt -eq "Import the CSV file provided on  OSF" -or $txt -match "^Import the CSV file provided on\s+OSF$") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$rng = $d.Range($target.Range.Start, $target.Range.End)

$runsXml = '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Import the </w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>Excel</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> file provided </w:t></w:r>' +
           '<w:proofErr w:type="gramStart"/>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">on  </w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>OSF</w:t></w:r>' +
           '<w:proofErr w:type="gramEnd"/>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>https://osf.io/jse8h</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>)</w:t></w:r>'

$xml = '<?xml version="1.0"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
